# Improve CT in Study Design sheet
#
# The "studyDesignBlindingScheme", "trialIntentTypes", "trialTypes" and
# "interventionModel" answers on the studyDesign sheet used to store raw
# controlled-terminology codes (e.g. "C49659=OPEN LABEL"). Replace them with
# the cleaned-up / more correct display text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("studyDesign")

$ws.Range("B3").Value = "OPEN LABEL"
$ws.Range("B4").Value = "BASIC SCIENCE,    DEVICE FEASIBILITY"
$ws.Range("B5").Value = "Efficacy Study"
$ws.Range("B6").Value = "C82639"

# Make studyDesign the active sheet/tab, with D13 selected (this also clears
# the previous tabSelected flag + selection on whichever sheet was active
# before).
$ws.Activate()
$ws.Range("D13").Select()
